$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the reported accuracy of row 5's measurement values from 3 to 2
# decimal places (custom accuracy).
$ws.Range("B5").Value = 1.92
$ws.Range("C5").Value = 1.11
$ws.Range("D5").Value = 0.62
$ws.Range("E5").Value = 3.89
$ws.Range("F5").Value = 3.27
$ws.Range("G5").Value = 1.52
$ws.Range("H5").Value = 14.22
$ws.Range("I5").Value = 2.33
$ws.Range("J5").Value = 1.04
$ws.Range("K5").Value = 1.43
$ws.Range("L5").Value = 1.49
$ws.Range("M5").Value = 1.59
$ws.Range("N5").Value = 0.51
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.23
$ws.Range("Q5").Value = 1.45
$ws.Range("R5").Value = 0.72
$ws.Range("S5").Value = 0.28
$ws.Range("T5").Value = 15.71
$ws.Range("U5").Value = 4.77
$ws.Range("V5").Value = 1.39
$ws.Range("W5").Value = 3.17
$ws.Range("X5").Value = 1.69
$ws.Range("Y5").Value = 0.18
$ws.Range("Z5").Value = 5.99
$ws.Range("AA5").Value = 1.23
$ws.Range("AB5").Value = 1.23
$ws.Range("AC5").Value = 1.41
$ws.Range("AD5").Value = 1.61
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 13.38
$ws.Range("AG5").Value = 0.7
$ws.Range("AH5").Value = 1.75

# Remove the last data row (row 6); Excel shifts the dimension ref down
# to A1:AH5 automatically.
$ws.Rows(6).Delete()
